# Applies the 2023-12-12 -> 2023-12-13 update to the "two digit division"
# worksheet: the date heading and every division problem in the table are
# replaced with the new day's values.

$d = $word.ActiveDocument

# --- Update the date heading ---------------------------------------------
$d.Content.Find.Execute("2023-12-12 Tuesday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2023-12-13 Wednesday", 2)

# --- Update the division problems in the table -----------------------------
# The worksheet table has 20 rows x 5 columns, but only every 4th row
# (1, 5, 9, 13, 17) actually holds a division problem - the rows in between
# are left blank for students to work in. Addressing cells directly by
# (row, column) avoids any ambiguity from repeated / reshuffled values.

$t = $d.Tables.Item(1)

$replacements = @{
    1  = @{ 1 = "71÷7="; 2 = "44÷2="; 3 = "25÷6="; 4 = "47÷7="; 5 = "58÷7=" }
    5  = @{ 1 = "32÷4="; 2 = "41÷5="; 3 = "47÷7="; 4 = "43÷6="; 5 = "59÷8=" }
    9  = @{ 1 = "14÷5="; 2 = "40÷4="; 3 = "39÷7="; 4 = "12÷5="; 5 = "19÷3=" }
    13 = @{ 1 = "70÷7="; 2 = "24÷3="; 3 = "28÷3="; 4 = "81÷2="; 5 = "67÷6=" }
    17 = @{ 1 = "42÷8="; 2 = "59÷7="; 3 = "94÷9="; 4 = "36÷9="; 5 = "88÷8=" }
}

foreach ($row in $replacements.Keys) {
    $cols = $replacements[$row]
    foreach ($col in $cols.Keys) {
        $t.Cell($row, $col).Range.Text = $cols[$col]
    }
}
